# Auto-generated Excel COM-interop script to apply market-price refresh values
# to the Gilgamesh_Profits workbook (columns H, I, J, K, L, M, N per Leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 422.57144
$ws.Range("I80").Value = 428.69232
$ws.Range("K80").Value = 1286.07696
$ws.Range("M80").Value = -288.0769599999999

$ws.Range("H83").Value = 422.57144
$ws.Range("I83").Value = 428.69232
$ws.Range("K83").Value = 3858.23088
$ws.Range("M83").Value = 1133.76912

$ws.Range("H88").Value = 20000858
$ws.Range("I88").Value = 50000348
$ws.Range("J88").Value = 1199
$ws.Range("K88").Value = 50000348
$ws.Range("L88").Value = 1199
$ws.Range("M88").Value = -49999942
$ws.Range("N88").Value = -2011

$ws.Range("H91").Value = 20000858
$ws.Range("I91").Value = 50000348
$ws.Range("J91").Value = 1199
$ws.Range("K91").Value = 50000348
$ws.Range("L91").Value = 1199
$ws.Range("M91").Value = -49998944
$ws.Range("N91").Value = -4007

$ws.Range("H98").Value = 2392.4119
$ws.Range("I98").Value = 3049.6924
$ws.Range("K98").Value = 3049.6924
$ws.Range("M98").Value = -1551.6924

$ws.Range("H122").Value = 2392.4119
$ws.Range("I122").Value = 3049.6924
$ws.Range("K122").Value = 9149.0772
$ws.Range("M122").Value = -6699.0772

$ws.Range("H137").Value = 911927.4
$ws.Range("I137").Value = 2382570.8
$ws.Range("J137").Value = 3588.7942
$ws.Range("K137").Value = 7147712.399999999
$ws.Range("L137").Value = 10766.3826
$ws.Range("M137").Value = -7145162.399999999
$ws.Range("N137").Value = -15866.3826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2923.3333
$ws.Range("I61").Value = 1125
$ws.Range("K61").Value = 1125
$ws.Range("M61").Value = -913

$ws.Range("H74").Value = 311384.78
$ws.Range("I74").Value = 926883.5
$ws.Range("J74").Value = 3635.4167
$ws.Range("K74").Value = 926883.5
$ws.Range("L74").Value = 3635.4167
$ws.Range("M74").Value = -926009.5
$ws.Range("N74").Value = -5383.4167

$ws.Range("H77").Value = 311384.78
$ws.Range("I77").Value = 926883.5
$ws.Range("J77").Value = 3635.4167
$ws.Range("K77").Value = 4634417.5
$ws.Range("L77").Value = 18177.0835
$ws.Range("M77").Value = -4630049.5
$ws.Range("N77").Value = -26913.0835

$ws.Range("H132").Value = 3577.4614
$ws.Range("I132").Value = 2627.5
$ws.Range("K132").Value = 7882.5
$ws.Range("M132").Value = -5352.5

$ws.Range("H135").Value = 114098
$ws.Range("J135").Value = 114098
$ws.Range("L135").Value = 114098
$ws.Range("N135").Value = -124238

$ws.Range("H136").Value = 2923.3333
$ws.Range("I136").Value = 1125
$ws.Range("K136").Value = 3375
$ws.Range("M136").Value = -825

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 74075256
$ws.Range("I94").Value = 100000360
$ws.Range("K94").Value = 100000360
$ws.Range("M94").Value = -99999909

$ws.Range("H107").Value = 1397.4584
$ws.Range("I107").Value = 1148.8823
$ws.Range("K107").Value = 1148.8823
$ws.Range("M107").Value = 771.1177

$ws.Range("H134").Value = 3420.4167
$ws.Range("I134").Value = 3287.3823
$ws.Range("K134").Value = 9862.1469
$ws.Range("M134").Value = -7327.1469

$ws.Range("H139").Value = 9536.846
$ws.Range("J139").Value = 9536.846
$ws.Range("L139").Value = 9536.846
$ws.Range("N139").Value = -19816.846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4167.7744
$ws.Range("J31").Value = 7400.8887
$ws.Range("L31").Value = 7400.8887
$ws.Range("N31").Value = -7990.8887

$ws.Range("H34").Value = 4167.7744
$ws.Range("J34").Value = 7400.8887
$ws.Range("L34").Value = 7400.8887
$ws.Range("N34").Value = -7804.8887

$ws.Range("H52").Value = 95606
$ws.Range("J52").Value = 95606
$ws.Range("L52").Value = 95606
$ws.Range("N52").Value = -96194

$ws.Range("H58").Value = 2029.0385
$ws.Range("I58").Value = 1035
$ws.Range("K58").Value = 1035
$ws.Range("M58").Value = -832

$ws.Range("H132").Value = 1598.4584
$ws.Range("I132").Value = 1348.7
$ws.Range("J132").Value = 2847.25
$ws.Range("K132").Value = 4046.1
$ws.Range("L132").Value = 8541.75
$ws.Range("M132").Value = -1516.1
$ws.Range("N132").Value = -13601.75

$ws.Range("H134").Value = 2175.4878
$ws.Range("I134").Value = 1951.5714
$ws.Range("J134").Value = 3481.6667
$ws.Range("K134").Value = 5854.7142
$ws.Range("L134").Value = 10445.0001
$ws.Range("M134").Value = -3319.7142
$ws.Range("N134").Value = -15515.0001

$ws.Range("H136").Value = 2029.0385
$ws.Range("I136").Value = 1035
$ws.Range("K136").Value = 3105
$ws.Range("M136").Value = -555

$ws.Range("H138").Value = 69998.836
$ws.Range("J138").Value = 69998.836
$ws.Range("L138").Value = 69998.836
$ws.Range("N138").Value = -80278.836

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 167557140
$ws.Range("I11").Value = 234400000
$ws.Range("K11").Value = 234400000
$ws.Range("M11").Value = -234399861

$ws.Range("H97").Value = 4850
$ws.Range("I97").Value = 9000
$ws.Range("J97").Value = 3466.6667
$ws.Range("K97").Value = 9000
$ws.Range("L97").Value = 3466.6667
$ws.Range("M97").Value = -8504
$ws.Range("N97").Value = -4458.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 7249.5
$ws.Range("I58").Value = 7249.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 7249.5
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("M58").Value = -6989.5

$ws.Range("H61").Value = 10279.615
$ws.Range("I61").Value = 2291
$ws.Range("K61").Value = 2291
$ws.Range("M61").Value = -2089

$ws.Range("H113").Value = 10279.615
$ws.Range("I113").Value = 2291
$ws.Range("K113").Value = 2291
$ws.Range("M113").Value = -121

$ws.Range("H133").Value = 73950.75
$ws.Range("J133").Value = 73950.75
$ws.Range("L133").Value = 73950.75
$ws.Range("N133").Value = -79010.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6325.375
$ws.Range("I132").Value = 6443.2856
$ws.Range("K132").Value = 19329.8568
$ws.Range("M132").Value = -16799.8568

$ws.Range("H138").Value = 115499
$ws.Range("J138").Value = 115499
$ws.Range("L138").Value = 115499
$ws.Range("N138").Value = -125779
